# Auto-generated edit script.
# Swaps the contents of row 2 and row 3 (all data columns, excluding formula
# cells which recompute automatically) across the nine per-fund compliance
# worksheets, so that fund HE3B2's figures move to row 2 and fund R21126's
# figures move to row 3 (previously the reverse).

$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Prospectus_80pct ---
$ws = $wb.Worksheets.Item('Prospectus_80pct')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'D2' 20531520.84
Set-CellValue $ws 'D3' 48822987.49999999
Set-CellValue $ws 'E2' 2080648.20256421
Set-CellValue $ws 'E3' 512.83423728
Set-CellValue $ws 'F2' -65189
Set-CellValue $ws 'F3' 0
Set-CellValue $ws 'G2' 0
Set-CellValue $ws 'G3' 0
Set-CellValue $ws 'H2' 45873.9990234375
Set-CellValue $ws 'H3' 103663.44921875
Set-CellValue $ws 'I2' 0
Set-CellValue $ws 'I3' 103150.61498147
Set-CellValue $ws 'J2' 22658043.04158765
Set-CellValue $ws 'J3' 48927163.78345602
Set-CellValue $ws 'K2' 20531520.84
Set-CellValue $ws 'K3' 48822987.49999999
Set-CellValue $ws 'M2' 20512205.83902344
Set-CellValue $ws 'M3' 48926650.94921874
Set-CellValue $ws 'N2' 20531520.84
Set-CellValue $ws 'N3' 48822987.49999999
Set-CellValue $ws 'P2' 'No'
Set-CellValue $ws 'P3' 'No'

# --- 40Act_Diversification ---
$ws = $wb.Worksheets.Item('40Act_Diversification')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 'diversified'
Set-CellValue $ws 'C3' 'diversified'
Set-CellValue $ws 'D2' 'diversified'
Set-CellValue $ws 'D3' 'diversified'
Set-CellValue $ws 'J2' 24293147.02
Set-CellValue $ws 'J3' 52846233.31
Set-CellValue $ws 'K2' 0
Set-CellValue $ws 'K3' 0
Set-CellValue $ws 'L2' 20717665.35
Set-CellValue $ws 'L3' 49252287.25
Set-CellValue $ws 'M2' 4491.06
Set-CellValue $ws 'M3' 10576.62
Set-CellValue $ws 'N2' 0
Set-CellValue $ws 'N3' 0
Set-CellValue $ws 'O2' 'None'
Set-CellValue $ws 'O3' 'None'
Set-CellValue $ws 'P2' 'NVDA, AAPL, MSFT, AVGO'
Set-CellValue $ws 'P3' 'GOOGLE, KLAC, AMAT, MLI, NVDA, LRCX, ROST, BK, BKR, AXP, EBAY'
Set-CellValue $ws 'Q2' 0.2616910394297878
Set-CellValue $ws 'Q3' 0.2659866179108261
Set-CellValue $ws 'R2' 0.7383089605702122
Set-CellValue $ws 'R3' 0.7340133820891739
Set-CellValue $ws 'S2' '(None, 0, 3.41%, 0.0000000%), (None, 0, 2.02%, 0.0000000%), (None, 0, 1.84%, 0.0000000%), (None, 0, 1.84%, 0.0000000%), (None, 0, 1.81%, 0.0000000%), (None, 0, 1.51%, 0.0000000%), (None, 0, 1.42%, 0.0000000%), (None, 0, 1.39%, 0.0000000%), (None, 0, 1.18%, 0.0000000%), (None, 0, 1.10%, 0.0000000%), (None, 0, 1.04%, 0.0000000%), (None, 0, 1.02%, 0.0000000%), (None, 0, 0.93%, 0.0000000%), (None, 0, 0.91%, 0.0000000%), (None, 0, 0.87%, 0.0000000%), (None, 0, 0.79%, 0.0000000%), (None, 0, 0.75%, 0.0000000%), (None, 0, 0.74%, 0.0000000%), (None, 0, 0.72%, 0.0000000%), (None, 0, 0.71%, 0.0000000%), (None, 0, 0.70%, 0.0000000%), (None, 0, 0.69%, 0.0000000%), (None, 0, 0.69%, 0.0000000%), (None, 0, 0.67%, 0.0000000%), (None, 0, 0.65%, 0.0000000%), (None, 0, 0.61%, 0.0000000%), (None, 0, 0.60%, 0.0000000%), (None, 0, 0.57%, 0.0000000%), (None, 0, 0.55%, 0.0000000%), (None, 0, 0.52%, 0.0000000%), (None, 0, 0.52%, 0.0000000%), (None, 0, 0.52%, 0.0000000%), (None, 0, 0.52%, 0.0000000%), (None, 0, 0.51%, 0.0000000%), (None, 0, 0.46%, 0.0000000%), (None, 0, 0.45%, 0.0000000%), (None, 0, 0.44%, 0.0000000%), (None, 0, 0.44%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.41%, 0.0000000%), (None, 0, 0.41%, 0.0000000%), (None, 0, 0.41%, 0.0000000%), (None, 0, 0.40%, 0.0000000%), (None, 0, 0.40%, 0.0000000%), (None, 0, 0.37%, 0.0000000%), (None, 0, 0.37%, 0.0000000%), (None, 0, 0.36%, 0.0000000%), (None, 0, 0.35%, 0.0000000%), (None, 0, 0.35%, 0.0000000%), (None, 0, 0.34%, 0.0000000%), (None, 0, 0.34%, 0.0000000%), (None, 0, 0.33%, 0.0000000%), (None, 0, 0.32%, 0.0000000%), (None, 0, 0.32%, 0.0000000%), (None, 0, 0.32%, 0.0000000%), (None, 0, 0.32%, 0.0000000%), (None, 0, 0.30%, 0.0000000%), (None, 0, 0.30%, 0.0000000%), (None, 0, 0.29%, 0.0000000%), (None, 0, 0.28%, 0.0000000%), (None, 0, 0.24%, 0.0000000%), (None, 0, 0.22%, 0.0000000%), (None, 0, 0.21%, 0.0000000%), (None, 0, 0.21%, 0.0000000%), (None, 0, 0.17%, 0.0000000%), (None, 0, 0.16%, 0.0000000%), (None, 0, 0.16%, 0.0000000%), (None, 0, 0.16%, 0.0000000%), (None, 0, 0.13%, 0.0000000%), (None, 0, 0.13%, 0.0000000%), (None, 0, 0.12%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%), (None, 0, nan%, 0.0000000%)'
Set-CellValue $ws 'S3' '(None, 0, 1.93%, 0.0000000%), (None, 0, 1.90%, 0.0000000%), (None, 0, 1.90%, 0.0000000%), (None, 0, 1.87%, 0.0000000%), (None, 0, 1.84%, 0.0000000%), (None, 0, 1.83%, 0.0000000%), (None, 0, 1.83%, 0.0000000%), (None, 0, 1.78%, 0.0000000%), (None, 0, 1.76%, 0.0000000%), (None, 0, 1.74%, 0.0000000%), (None, 0, 1.73%, 0.0000000%), (None, 0, 1.73%, 0.0000000%), (None, 0, 1.69%, 0.0000000%), (None, 0, 1.67%, 0.0000000%), (None, 0, 1.66%, 0.0000000%), (None, 0, 1.64%, 0.0000000%), (None, 0, 1.59%, 0.0000000%), (None, 0, 1.57%, 0.0000000%), (None, 0, 1.56%, 0.0000000%), (None, 0, 1.53%, 0.0000000%), (None, 0, 1.47%, 0.0000000%), (None, 0, 1.45%, 0.0000000%), (None, 0, 1.42%, 0.0000000%), (None, 0, 1.40%, 0.0000000%), (None, 0, 1.40%, 0.0000000%), (None, 0, 1.37%, 0.0000000%), (None, 0, 1.30%, 0.0000000%), (None, 0, 1.28%, 0.0000000%), (None, 0, 1.24%, 0.0000000%), (None, 0, 1.24%, 0.0000000%), (None, 0, 1.04%, 0.0000000%), (None, 0, 0.91%, 0.0000000%), (None, 0, 0.91%, 0.0000000%), (None, 0, 0.88%, 0.0000000%), (None, 0, 0.87%, 0.0000000%), (None, 0, 0.87%, 0.0000000%), (None, 0, 0.81%, 0.0000000%), (None, 0, 0.78%, 0.0000000%), (None, 0, 0.76%, 0.0000000%), (None, 0, 0.75%, 0.0000000%), (None, 0, 0.68%, 0.0000000%), (None, 0, 0.65%, 0.0000000%), (None, 0, 0.62%, 0.0000000%), (None, 0, 0.57%, 0.0000000%), (None, 0, 0.55%, 0.0000000%), (None, 0, 0.51%, 0.0000000%), (None, 0, 0.51%, 0.0000000%), (None, 0, 0.50%, 0.0000000%), (None, 0, 0.47%, 0.0000000%), (None, 0, 0.46%, 0.0000000%), (None, 0, 0.44%, 0.0000000%), (None, 0, 0.44%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.43%, 0.0000000%), (None, 0, 0.42%, 0.0000000%), (None, 0, 0.42%, 0.0000000%), (None, 0, 0.39%, 0.0000000%), (None, 0, 0.36%, 0.0000000%), (None, 0, 0.36%, 0.0000000%), (None, 0, 0.35%, 0.0000000%), (None, 0, 0.30%, 0.0000000%)'
Set-CellValue $ws 'T2' 0
Set-CellValue $ws 'T3' 0
Set-CellValue $ws 'U2' 65189
Set-CellValue $ws 'U3' 0
Set-CellValue $ws 'V2' 0.002683431666812512
Set-CellValue $ws 'V3' 0

# --- IRS_Diversification ---
$ws = $wb.Worksheets.Item('IRS_Diversification')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 'PASS'
Set-CellValue $ws 'C3' 'PASS'
Set-CellValue $ws 'F2' 'PASS'
Set-CellValue $ws 'F3' 'PASS'
Set-CellValue $ws 'G2' 'PASS'
Set-CellValue $ws 'G3' 'PASS'
Set-CellValue $ws 'H2' 24293147.02
Set-CellValue $ws 'H3' 52846233.31
Set-CellValue $ws 'I2' 16651926.38
Set-CellValue $ws 'I3' 0
Set-CellValue $ws 'J2' 4491.06
Set-CellValue $ws 'J3' 10576.62
Set-CellValue $ws 'K2' 1214657.351
Set-CellValue $ws 'K3' 2642311.6655
Set-CellValue $ws 'L2' 0.2210311491492452
Set-CellValue $ws 'L3' 0
Set-CellValue $ws 'M2' 3
Set-CellValue $ws 'M3' 0
Set-CellValue $ws 'N2' '(, 8.12%), (, 7.53%), (, 6.46%)'
Set-CellValue $ws 'N3' 'None'
Set-CellValue $ws 'O2' 0
Set-CellValue $ws 'O3' 0
Set-CellValue $ws 'P2' 'None'
Set-CellValue $ws 'P3' 'None'
Set-CellValue $ws 'Q2' 'ABBV (1.09%)'
Set-CellValue $ws 'Q3' 'N/A (0.00%)'
Set-CellValue $ws 'R2' 'NEE (1.07%)'
Set-CellValue $ws 'R3' 'N/A (0.00%)'

# --- Illiquid ---
$ws = $wb.Worksheets.Item('Illiquid')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 24293147.02
Set-CellValue $ws 'C3' 52846233.31
Set-CellValue $ws 'D2' 0
Set-CellValue $ws 'D3' 0
Set-CellValue $ws 'E2' 0
Set-CellValue $ws 'E3' 0
Set-CellValue $ws 'F2' 0.8451569005488201
Set-CellValue $ws 'F3' 0.9238688254960511

# --- Real_Estate ---
$ws = $wb.Worksheets.Item('Real_Estate')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 0
Set-CellValue $ws 'C3' 0
Set-CellValue $ws 'D2' 20531520.84
Set-CellValue $ws 'D3' 48822987.49999999

# --- Commodities ---
$ws = $wb.Worksheets.Item('Commodities')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 0
Set-CellValue $ws 'C3' 0

# --- 12d1_Other_Investment_Companies ---
$ws = $wb.Worksheets.Item('12d1_Other_Investment_Companies')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 24293147.02
Set-CellValue $ws 'C3' 52846233.31
Set-CellValue $ws 'D2' 'None'
Set-CellValue $ws 'D3' 'None'
Set-CellValue $ws 'E2' 0
Set-CellValue $ws 'E3' 0
Set-CellValue $ws 'F2' 0
Set-CellValue $ws 'F3' 0

# --- 12d2_Insurance_Companies ---
$ws = $wb.Worksheets.Item('12d2_Insurance_Companies')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'C2' 'PASS'
Set-CellValue $ws 'C3' 'PASS'
Set-CellValue $ws 'D2' 24293147.02
Set-CellValue $ws 'D3' 52846233.31
Set-CellValue $ws 'E2' 'CB (0.00000%), MET (0.00000%)'
Set-CellValue $ws 'E3' 'AFL (0.00000%), ALL (0.00000%), CB (0.00000%), CINF (0.00000%), EG (0.00000%), HIG (0.00000%), LNC (0.00000%), MET (0.00000%), PGR (0.00000%), TRV (0.00000%), UNM (0.00000%)'

# --- 12d3_Securities_Business ---
$ws = $wb.Worksheets.Item('12d3_Securities_Business')
Set-CellValue $ws 'B2' 'HE3B2'
Set-CellValue $ws 'B3' 'R21126'
Set-CellValue $ws 'D2' 'PASS'
Set-CellValue $ws 'D3' 'PASS'
Set-CellValue $ws 'G2' 'None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%)'
Set-CellValue $ws 'G3' 'None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%), None (0.00000%)'
Set-CellValue $ws 'H2' 'None (0.75%), None (0.39%), None (0.41%), None (0.65%), None (1.84%), None (0.52%), None (0.36%), None (0.34%), None (0.16%), None (0.60%)'
Set-CellValue $ws 'H3' 'None (0.75%), None (0.47%), None (2.06%), None (1.30%), None (1.84%), None (1.93%), None (0.42%), None (0.51%), None (1.28%), None (1.67%), None (1.45%), None (0.42%)'
Set-CellValue $ws 'I2' 24293147.02
Set-CellValue $ws 'I3' 52846233.31
Set-CellValue $ws 'J2' 0
Set-CellValue $ws 'J3' 0
Set-CellValue $ws 'K2' 0.01835160095285177
Set-CellValue $ws 'K3' 0.02058899398974023

